$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted above the old row 22, pushing the old
# row 22 data down to row 23 (dimension grows from A1:R22 to A1:R23).
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the latest weekly observation.
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value = "Ñuble"
$ws.Range("D22").Value = 44509
$ws.Range("D22").NumberFormat = $ws.Range("D21").NumberFormat
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 100112022
$ws.Range("G22").Value = "Arveja Verde"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 17500
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 700
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
